$d = $word.ActiveDocument

# Using Find (match-only, action-less) + direct Range.Text assignment rather
# than Find.Execute's Replace parameter avoids Word's smart-quote
# autocorrect kicking in on straight apostrophes in the replacement text.

function Replace-ExactText($oldText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = $newText
    }
    return $found
}

$old0 = "I am a senior AI/ML Engineer with a decade of experience in the software development industry, specializing in AI, MLOps, and system architecture. My journey in AI began in America, where I have consistently been at the forefront of technological advancements, architecting and implementing cutting-edge software solutions tailored to business needs. My expertise in AI-driven platforms and projects has enabled me to streamline business operations, enhance customer service, and drive innovation."
$new0 = "I am a senior AI/ML Engineer with extensive experience in building scalable GenAI platforms, cloud-native ML pipelines, and high-impact AI solutions. With over eight years in the field, I have honed my skills in AWS and Kubernetes, making me an ideal candidate for the Data Platform Infrastructure Engineer position. My expertise in orchestrating CI/CD pipelines, automated testing frameworks, and real-time monitoring systems aligns perfectly with your requirement for ensuring the reliability, scalability, and performance of large-scale systems."

$old1 = "At InsoftAI, I led the development of AI-driven platforms capable of handling up to 90% of customer inquiries, significantly enhancing operational efficiency. My work on Support-nGen" + [char]0x2122 + " and LLM Twin showcases my ability to develop proprietary systems that automate processes and generate creative ideas, facilitating brand creation and streamlining content creation. My proficiency in deploying scalable, secure, and efficient real-time predictions on AWS SageMaker further highlights my technical capabilities."
$new1 = "My recent role at Fingent involved architecting CI/CD pipelines using Jenkins and GitLab CI/CD, which improved deployment times by 30% for a multinational logistics client. I have also optimized cloud infrastructure with Terraform on AWS, resulting in a 20% reduction in cloud costs while accommodating a 50% workload increase. My hands-on experience with AWS's ecosystem, including SageMaker and ECR, has equipped me to effectively manage and scale infrastructure."

$old2 = "During my tenure at Brainhub, I developed Sierra.ai, revolutionizing document management and information accessibility, resulting in a 30% increase in operational efficiency for clients. My leadership in formulating technical strategies and optimizing multi-AI agents improved response times by 40% and reduced manual intervention. I have a proven track record of designing and implementing robust ML serving architectures and deploying scalable, cost-effective solutions that align with business goals."
$new2 = "I thrive in agile environments and excel in collaborating with cross-functional teams to deliver robust data platforms. My programming skills in Python and Java, coupled with my ability to adapt to changes, make me a fast learner and a valuable asset to your team. Additionally, my experience in MLOps and machine learning pipelines further enhances my capability to contribute effectively to your organization's goals."

$old3 = "My experience at Kensho involved building TTS and STT solutions, enhancing user experience in voice synthesis applications, and developing ML systems for forecasting energy consumption. I have demonstrated strong leadership abilities by mentoring junior staff and fostering skill development, enhancing team performance."
$new3 = "I am eager to bring my expertise in AWS, Kubernetes, and infrastructure-as-code to your team, ensuring the development of cutting-edge tools and dashboards for monitoring and management. Thank you for considering my application, and I look forward to the opportunity to contribute to your innovative projects."

$old4 = "I am eager to define and drive the long-term ML technical strategy in alignment with product and business goals. My deep expertise in computer vision, GenAI, and adjacent fields, combined with my ability to lead and grow high-performing teams, makes me an ideal candidate for this role. I am committed to creating a team culture where people feel empowered, supported, and technically challenged, ensuring strong cross-functional collaboration and delivering state-of-the-art models into production swiftly. I look forward to contributing to your organization's success by leveraging my skills and experience in AI/ML engineering."

Replace-ExactText $old0 $new0 | Out-Null
Replace-ExactText $old1 $new1 | Out-Null
Replace-ExactText $old2 $new2 | Out-Null
Replace-ExactText $old3 $new3 | Out-Null

# --- Remove the trailing (now orphaned) 5th paragraph block, together with
#     the two <w:br/> line breaks that introduced it, so the new 4th block
#     becomes the final block of the run. ---

$rng = $d.Content
$found = $rng.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $delStart = $rng.Start - 2
    $delEnd = $rng.End
    $delRange = $d.Range($delStart, $delEnd)
    $delRange.Text = ""
}
